# The presentation's active Design ("Integral" theme / "Red Violet" colour
# scheme - stored in ppt/theme/theme2.xml, the theme actually wired up to
# the slide master and to the presentation itself) is switched back to the
# default "Office Theme" / "Office" colour scheme. The font scheme and the
# fill/line/effect format scheme are already identical between the themes
# in this deck, so only the twelve theme colours need to change.

function ColorRef([int]$r, [int]$g, [int]$b) {
    # PowerPoint's ColorFormat.RGB uses the OLE COLORREF packing (0x00BBGGRR)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Index order matches the OOXML <a:clrScheme> child order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$colorScheme.Item(1).RGB  = ColorRef 0x00 0x00 0x00   # dk1      000000
$colorScheme.Item(2).RGB  = ColorRef 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = ColorRef 0x44 0x54 0x6A   # dk2      44546A
$colorScheme.Item(4).RGB  = ColorRef 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = ColorRef 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = ColorRef 0xED 0x7D 0x31   # accent2  ED7D31
$colorScheme.Item(7).RGB  = ColorRef 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = ColorRef 0xFF 0xC0 0x00   # accent4  FFC000
$colorScheme.Item(9).RGB  = ColorRef 0x44 0x72 0xC4   # accent5  4472C4
$colorScheme.Item(10).RGB = ColorRef 0x70 0xAD 0x47   # accent6  70AD47
$colorScheme.Item(11).RGB = ColorRef 0x05 0x63 0xC1   # hlink    0563C1
$colorScheme.Item(12).RGB = ColorRef 0x95 0x4F 0x72   # folHlink 954F72
